# Cambio y correccion en numero de cuentas por cobrar cliente
#
# 1. Rename the "Generar Adeudo" label (I2) to make clear it also covers
#    "Cuentas por Cobrar Clientes".
# 2. Mark the "Alicuotas" row (row 15) as an accumulating account (D15 = "x"),
#    same convention used for the other accumulating accounts (D6, D8, D9).
# 3. Fix/extend the "Cuentas por Cobrar" detail accounts: instead of the old,
#    incorrect single "CxC condominos" (1.1.3.01) line, list the accounts per
#    departamento: Departamento 101 (1.1.3.01) and the new Departamento 102
#    (1.1.3.02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename "Generar Adeudo" -> "Generar Adeudo / Ctas.por Cobrar Clientes"
$ws.Range("I2").Value = "Generar Adeudo / Ctas.por Cobrar Clientes"

# --- 2. Flag Alicuotas (row 15) as an "Acumula" account, like D6/D8/D9
$ws.Range("D15").Value = "x"
$ws.Range("D15").HorizontalAlignment = -4108  # xlCenter, matches D6/D8/D9 style

# --- 3. Rebuild rows 16-17 with the corrected Cuentas por Cobrar detail
$ws.Range("A16").Value = "Departamento 101"
$ws.Range("B16").Value = "1.1.3.01"
$ws.Range("C16").Value = 4

$ws.Range("B17").Value = "1.1.3.02"
$ws.Range("A17").Value = "Departamento 102"
$ws.Range("C17").Value = 4

$titleRange = $ws.Range("A16:A17")
$titleRange.Font.Bold = $false
$titleRange.HorizontalAlignment = -4108  # xlCenter

# --- cosmetic: leave selection where the author left it
$ws.Range("D15").Select() | Out-Null
